$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 11:38"

# Update country labels that shifted rank due to refreshed data
$ws.Range("A110").Value = "Albania"
$ws.Range("A111").Value = "Eslovaquia"
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"

# Update numeric statistics for affected rows
# Row 7
$ws.Range("B7").Value = 333255
$ws.Range("C7").Value = 472
$ws.Range("D7").Value = 169817
$ws.Range("E7").Value = 153914
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 9524

# Row 21
$ws.Range("B21").Value = 90619
$ws.Range("C21").Value = 3099
$ws.Range("E21").Value = 70680
$ws.Range("G21").Value = 38
$ws.Range("H21").Value = 1209

# Row 25
$ws.Range("B25").Value = 60100
$ws.Range("C25").Value = 71
$ws.Range("D25").Value = 16610
$ws.Range("E25").Value = 33829
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 9661

# Row 34
$ws.Range("B34").Value = 39294
$ws.Range("C34").Value = 1017
$ws.Range("D34").Value = 15123
$ws.Range("E34").Value = 21973
$ws.Range("G34").Value = 64
$ws.Range("H34").Value = 2198

# Row 40
$ws.Range("B40").Value = 29788
$ws.Range("C40").Value = 396
$ws.Range("E40").Value = 14149
$ws.Range("G40").Value = 9
$ws.Range("H40").Value = 1256

# Row 41
$ws.Range("B41").Value = 26420
$ws.Range("C41").Value = 490
$ws.Range("D41").Value = 6252
$ws.Range("E41").Value = 19070
$ws.Range("G41").Value = 10
$ws.Range("H41").Value = 1098

# Row 44
$ws.Range("B44").Value = 24524
$ws.Range("C44").Value = 1043
$ws.Range("D44").Value = 9533
$ws.Range("E44").Value = 14883
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 108

# Row 47
$ws.Range("B47").Value = 20686
$ws.Range("E47").Value = 6483

# Row 53
$ws.Range("B53").Value = 17135
$ws.Range("C53").Value = 26
$ws.Range("D53").Value = 16066
$ws.Range("E53").Value = 391
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 678

# Row 63
$ws.Range("B63").Value = 10027
$ws.Range("C63").Value = 3
$ws.Range("D63").Value = 7245
$ws.Range("E63").Value = 2452

# Row 68
$ws.Range("B68").Value = 8838
$ws.Range("C68").Value = 45
$ws.Range("D68").Value = 7779
$ws.Range("E68").Value = 847

# Row 70
$ws.Range("B70").Value = 8494
$ws.Range("C70").Value = 41
$ws.Range("D70").Value = 7400
$ws.Range("E70").Value = 973

# Row 86
$ws.Range("B86").Value = 3826
$ws.Range("C86").Value = 106
$ws.Range("D86").Value = 1912
$ws.Range("E86").Value = 1840

# Row 103
$ws.Range("D103").Value = 1342
$ws.Range("E103").Value = 536

# Row 110
$ws.Range("B110").Value = 1590
$ws.Range("C110").Value = 69
$ws.Range("D110").Value = 1055
$ws.Range("E110").Value = 499
$ws.Range("H110").Value = 36

# Row 111
$ws.Range("B111").Value = 1552
$ws.Range("C111").Value = 4
$ws.Range("D111").Value = 1410
$ws.Range("E111").Value = 114
$ws.Range("H111").Value = 28

# Row 138
$ws.Range("B138").Value = 705
$ws.Range("C138").Value = 9
$ws.Range("D138").Value = 299
$ws.Range("E138").Value = 406

# Row 208
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
